# Weekly CompStat update: new crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
# "Volume 32   Number  22" -> "...23"
$ws.Range("C1").Value = "Volume 32   Number  23"

# "Report Covering the Week  5/26/2025  Through  6/1/2025"
#   -> "...6/2/2025  Through  6/8/2025"
$ws.Range("A6").Value = "Report Covering the Week  6/2/2025  Through  6/8/2025"

# --- Row 16 (Murder) ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = -26.470588235294
$ws.Range("L16").Value = -10.714285714285
$ws.Range("M16").Value = -24.242424242424
$ws.Range("N16").Value = -82.758620689655

# --- Row 17 (Rape) ---
# C17 switches from numeric 1 to text "0" (same style family as other
# text placeholder cells, e.g. A17).
$ws.Range("C17").Value = "'0"
$ws.Range("A17").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -58.333333333333
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = -21.276595744680
$ws.Range("L17").Value = -22.916666666666
$ws.Range("M17").Value = -11.904761904761
$ws.Range("N17").Value = -71.317829457364

# --- Row 18 (Robbery) ---
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -44.444444444444
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = 17.647058823529
$ws.Range("L18").Value = 2.564102564102
$ws.Range("M18").Value = -11.111111111111
$ws.Range("N18").Value = -77.900552486187

# --- Row 19 (Fel. Assault) ---
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 17
$ws.Range("H19").Value = -10.526315789473
$ws.Range("I19").Value = 80
$ws.Range("J19").Value = 87
$ws.Range("K19").Value = -8.045977011494
$ws.Range("L19").Value = 17.647058823529
$ws.Range("M19").Value = -4.761904761904
$ws.Range("N19").Value = -6.976744186046

# --- Row 20 (Burglary) ---
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -50
$ws.Range("J20").Value = 21
$ws.Range("K20").Value = -4.761904761904
$ws.Range("L20").Value = -35.483870967741
$ws.Range("M20").Value = -28.571428571428
$ws.Range("N20").Value = -89.071038251366

# --- Row 21 (Gr. Larceny) ---
$ws.Range("C21").Value = 5
$ws.Range("E21").Value = -66.666666666666
$ws.Range("F21").Value = 37
$ws.Range("G21").Value = 51
$ws.Range("H21").Value = -27.450980392156
$ws.Range("I21").Value = 203
$ws.Range("J21").Value = 226
$ws.Range("K21").Value = -10.176991150442
$ws.Range("L21").Value = -6.018518518518
$ws.Range("M21").Value = -13.617021276595
$ws.Range("N21").Value = -72.305593451568

# --- Row 23 (TOTAL) ---
# C23 switches from numeric 2 to text "0".
$ws.Range("C23").Value = "'0"
$ws.Range("A23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -37.5
$ws.Range("J23").Value = 47
$ws.Range("K23").Value = -14.893617021276
$ws.Range("L23").Value = -20
$ws.Range("M23").Value = 42.857142857142

# --- Row 24 (Transit) ---
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -76
$ws.Range("F24").Value = 33
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = -54.166666666666
$ws.Range("I24").Value = 206
$ws.Range("J24").Value = 226
$ws.Range("K24").Value = -8.849557522123
$ws.Range("L24").Value = -23.134328358209
$ws.Range("M24").Value = 5.102040816326

# --- Row 25 (Housing) ---
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -91.304347826087
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = -81.034482758620
$ws.Range("I25").Value = 96
$ws.Range("J25").Value = 126
$ws.Range("K25").Value = -23.809523809523
$ws.Range("L25").Value = -35.570469798657

# --- Row 26 (Petit Larceny) ---
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 7.692307692307
$ws.Range("I26").Value = 81
$ws.Range("J26").Value = 62
$ws.Range("K26").Value = 30.645161290322
$ws.Range("L26").Value = 17.391304347826
$ws.Range("M26").Value = -33.606557377049

# --- Row 27 (Retail Theft) ---
$ws.Range("L27").Value = -66.666666666666

# --- Row 28 (Misd. Assault) ---
# C28 and D28 switch from numeric to text "0"; E28 switches from numeric
# -50 to text "***.*" -- all three take on the same "text placeholder"
# style as A28. F28/G28/H28 are plain numeric updates.
$ws.Range("C28").Value = "'0"
$ws.Range("D28").Value = "'0"
$ws.Range("E28").Value = "***.*"
$ws.Range("A28").Copy()
$ws.Range("C28:E28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75

# --- Row 29 (UCR Rape*) ---
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -85

# --- Row 30 (Other Sex Crimes) ---
$ws.Range("M30").Value = -66.666666666666
$ws.Range("N30").Value = -85.714285714285
